# Hortaliza, Macroferia Regional de Talca - Berenjena
# Inserts one new weekly price-report row (new row 164) above the former
# row 164, shifting the existing rows 164:186 down to 165:187.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 164 (pushes old 164..186 down to 165..187)
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with the new record
$ws.Range("A164").Value = 5
$ws.Range("B164").Value = "Macroferia Regional de Talca"
$ws.Range("C164").Value = "Maule"
$ws.Range("D164").Value = 45131
$ws.Range("E164").Value = 7
$ws.Range("F164").Value = 100112001
$ws.Range("G164").Value = "Berenjena"
$ws.Range("H164").Value = "Sin especificar"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 200
$ws.Range("K164").Value = 8000
$ws.Range("L164").Value = 8000
$ws.Range("M164").Value = 8000
$ws.Range("N164").Value = "`$/caja 50 unidades"
$ws.Range("O164").Value = "Región de Arica y Parinacota"
$ws.Range("P164").Value = 160
$ws.Range("Q164").Value = 50
$ws.Range("R164").Value = "Hortaliza"
